$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with revised values (675-677)
$ws.Range("F675").Value = 3798.03799633
$ws.Range("E676").Value = 29190.55
$ws.Range("F676").Value = 13725.97767151
$ws.Range("F677").Value = 22921.12890833

# Style template cell for column A (date) formatting used by existing rows
$styleSrc = $ws.Range("A677")

# Append new rows 678-691
$styleSrc.Copy($ws.Range("A678"))
$ws.Range("A678").Value = 45147.41666666666
$ws.Range("B678").Value = 29775.65
$ws.Range("C678").Value = 30129.27
$ws.Range("D678").Value = 29365.49
$ws.Range("E678").Value = 29573.89
$ws.Range("F678").Value = 15359.03749689

$styleSrc.Copy($ws.Range("A679"))
$ws.Range("A679").Value = 45148.41666666666
$ws.Range("B679").Value = 29573.92
$ws.Range("C679").Value = 29712.95
$ws.Range("D679").Value = 29317.25
$ws.Range("E679").Value = 29433.51
$ws.Range("F679").Value = 11115.30974098

$styleSrc.Copy($ws.Range("A680"))
$ws.Range("A680").Value = 45149.41666666666
$ws.Range("B680").Value = 29434.01
$ws.Range("C680").Value = 29537.54
$ws.Range("D680").Value = 29223.42
$ws.Range("E680").Value = 29407.86
$ws.Range("F680").Value = 8933.77670339

$styleSrc.Copy($ws.Range("A681"))
$ws.Range("A681").Value = 45150.41666666666
$ws.Range("B681").Value = 29407.86
$ws.Range("C681").Value = 29473.73
$ws.Range("D681").Value = 29361.72
$ws.Range("E681").Value = 29422.34
$ws.Range("F681").Value = 2842.16799777

$styleSrc.Copy($ws.Range("A682"))
$ws.Range("A682").Value = 45151.41666666666
$ws.Range("B682").Value = 29422.42
$ws.Range("C682").Value = 29451.93
$ws.Range("D682").Value = 29264.29
$ws.Range("E682").Value = 29289.76
$ws.Range("F682").Value = 3248.0312529

$styleSrc.Copy($ws.Range("A683"))
$ws.Range("A683").Value = 45152.41666666666
$ws.Range("B683").Value = 29288.97
$ws.Range("C683").Value = 29667.77
$ws.Range("D683").Value = 29090.49
$ws.Range("E683").Value = 29419.22
$ws.Range("F683").Value = 10383.36599966

$styleSrc.Copy($ws.Range("A684"))
$ws.Range("A684").Value = 45153.41666666666
$ws.Range("B684").Value = 29419
$ws.Range("C684").Value = 29467.16
$ws.Range("D684").Value = 29064.65
$ws.Range("E684").Value = 29176.89
$ws.Range("F684").Value = 8584.02071459

$styleSrc.Copy($ws.Range("A685"))
$ws.Range("A685").Value = 45154.41666666666
$ws.Range("B685").Value = 29176.63
$ws.Range("C685").Value = 29232.71
$ws.Range("D685").Value = 28701.67
$ws.Range("E685").Value = 28707.5
$ws.Range("F685").Value = 14000.5675284

$styleSrc.Copy($ws.Range("A686"))
$ws.Range("A686").Value = 45155.41666666666
$ws.Range("B686").Value = 28705.2
$ws.Range("C686").Value = 28758.96
$ws.Range("D686").Value = 25253.44
$ws.Range("E686").Value = 26635.04
$ws.Range("F686").Value = 43569.94033756

$styleSrc.Copy($ws.Range("A687"))
$ws.Range("A687").Value = 45156.41666666666
$ws.Range("B687").Value = 26631.58
$ws.Range("C687").Value = 26824.09
$ws.Range("D687").Value = 25618.28
$ws.Range("E687").Value = 26053.12
$ws.Range("F687").Value = 28983.00011015

$styleSrc.Copy($ws.Range("A688"))
$ws.Range("A688").Value = 45157.41666666666
$ws.Range("B688").Value = 26054.35
$ws.Range("C688").Value = 26267.78
$ws.Range("D688").Value = 25800.8
$ws.Range("E688").Value = 26097.91
$ws.Range("F688").Value = 8854.32722316

$styleSrc.Copy($ws.Range("A689"))
$ws.Range("A689").Value = 45158.41666666666
$ws.Range("B689").Value = 26096.9
$ws.Range("C689").Value = 26295.77
$ws.Range("D689").Value = 25987.68
$ws.Range("E689").Value = 26196.16
$ws.Range("F689").Value = 6240.6993679

$styleSrc.Copy($ws.Range("A690"))
$ws.Range("A690").Value = 45159.41666666666
$ws.Range("B690").Value = 26195.97
$ws.Range("C690").Value = 26251.06
$ws.Range("D690").Value = 25820.83
$ws.Range("E690").Value = 26129.39
$ws.Range("F690").Value = 13690.79828458

$styleSrc.Copy($ws.Range("A691"))
$ws.Range("A691").Value = 45160.41666666666
$ws.Range("B691").Value = 26129.39
$ws.Range("C691").Value = 26138.6
$ws.Range("D691").Value = 25361.73
$ws.Range("E691").Value = 26046.38
$ws.Range("F691").Value = 16916.00830198
